# Apply updated crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe
$ws.Range("D2").Value = '68.694.11'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '2.444.87'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("E6").Value = '  +2.17%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +3.05%  '
$ws.Range("E9").Value = '  +8.77%  '
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("E11").Value = '  -1.41%  '
$ws.Range("E12").Value = '  +2.14%  '
$ws.Range("D13").Value = '68.572.45'
$ws.Range("E13").Value = '  +1.03%  '
$ws.Range("E14").Value = '  +3.93%  '
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("E16").Value = '  -2.43%  '
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("E19").Value = '  +2.09%  '
$ws.Range("B20").Value = 'SuiNetwork'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("E20").Value = '  +3.15%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +1.45%  '
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("E24").Value = '  +2.28%  '
$ws.Range("D25").Value = '0.0₃0816'
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("E26").Value = '  +2.19%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("E29").Value = '  +1.98%  '
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("E37").Value = '  +1.73%  '
$ws.Range("E38").Value = '  +2.70%  '
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("E40").Value = '  +1.27%  '
$ws.Range("E41").Value = '  +2.26%  '
$ws.Range("E42").Value = '  -1.87%  '
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("E44").Value = '  +1.57%  '
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("D51").Value = '0.0₆0204'
$ws.Range("E51").Value = '  +5.09%  '

# Numeric-looking values that must stay stored as text (preserve trailing
# zeros / exact formatting, matching the original inline-string cells):
# temporarily force a Text number format, assign, then restore the cell
# style so no stray formatting sticks to the cell.
$textCells = @{
    'D4' = '0.999'
    'D5' = '559.90'
    'D6' = '163.42'
    'D8' = '0.512'
    'D9' = '0.158'
    'D11' = '0.328'
    'D12' = '4.80'
    'D14' = '0.0000171'
    'D15' = '23.28'
    'D16' = '10.41'
    'D17' = '337.64'
    'D18' = '6.89'
    'D19' = '3.80'
    'D20' = '1.89'
    'D21' = '1.00'
    'D22' = '66.97'
    'D23' = '3.69'
    'D24' = '8.15'
    'D26' = '7.20'
    'D28' = '426.35'
    'D31' = '161.04'
    'D32' = '19.01'
    'D34' = '17.80'
    'D36' = '0.297'
    'D37' = '4.37'
    'D40' = '2.02'
    'D41' = '3.36'
    'D42' = '129.86'
    'D43' = '0.0716'
    'D44' = '0.480'
    'D45' = '0.562'
    'D46' = '0.0919'
    'D49' = '4.91'
    'D50' = '16.71'
}
foreach ($cellref in $textCells.Keys) {
    $r = $ws.Range($cellref)
    $r.NumberFormat = "@"
    $r.Value = $textCells[$cellref]
    $r.Style = "Normal"
}
